$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(3,4,6,7,8,9,11,12,15,16,17,18,25,28,29,30,31)

foreach ($r in $rows) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq "Sown Permanent Pasture") {
        $cell.Value = "Natural Pasture"
    }
}
